# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
#
# The worker-account table (rows 16-20, columns B:G) is refreshed:
# CLELIA ELENA PUELLO DIAZ is moved to the top of the list (row 16) and
# the remaining workers shift down one position; the "Valor Mora"
# (F) and "Salario Basico" (G) amounts are updated to their new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: CLELIA ELENA PUELLO DIAZ
$ws.Cells.Item(16, 2).Value = "CC"
$ws.Cells.Item(16, 3).Value = "45591669"
$ws.Cells.Item(16, 4).Value = "CLELIA ELENA PUELLO DIAZ"
$ws.Cells.Item(16, 5).Value = "1704"
$ws.Cells.Item(16, 6).Value = 140000
$ws.Cells.Item(16, 7).Value = 4150000

# Row 17: SALOMON DE JESUS DAGUER QUINTERO
$ws.Cells.Item(17, 2).Value = "CC"
$ws.Cells.Item(17, 3).Value = "73120289"
$ws.Cells.Item(17, 4).Value = "SALOMON DE JESUS DAGUER QUINTERO"
$ws.Cells.Item(17, 5).Value = "1704"
$ws.Cells.Item(17, 6).Value = 220000
$ws.Cells.Item(17, 7).Value = 5500000

# Row 18: RAFAEL GUILLERMO DAGUER QUINTERO
$ws.Cells.Item(18, 2).Value = "CC"
$ws.Cells.Item(18, 3).Value = "73101373"
$ws.Cells.Item(18, 4).Value = "RAFAEL GUILLERMO DAGUER QUINTERO"
$ws.Cells.Item(18, 5).Value = "1704"
$ws.Cells.Item(18, 6).Value = 360000
$ws.Cells.Item(18, 7).Value = 9000000

# Row 19: DINOHORA JUDITH CHAVEZ MORA
$ws.Cells.Item(19, 2).Value = "CC"
$ws.Cells.Item(19, 3).Value = "33103945"
$ws.Cells.Item(19, 4).Value = "DINOHORA JUDITH CHAVEZ MORA"
$ws.Cells.Item(19, 5).Value = "1704"
$ws.Cells.Item(19, 6).Value = 128000
$ws.Cells.Item(19, 7).Value = 3200000

# Row 20: JUAN CARLOS VARELA SIERRA
$ws.Cells.Item(20, 2).Value = "CC"
$ws.Cells.Item(20, 3).Value = "7919969"
$ws.Cells.Item(20, 4).Value = "JUAN CARLOS VARELA SIERRA"
$ws.Cells.Item(20, 5).Value = "1704"
$ws.Cells.Item(20, 6).Value = 220000
$ws.Cells.Item(20, 7).Value = 5500000
